# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-48, columns E:G) was re-sorted from
# descending chronological order (2104 .. 1704) to ascending chronological
# order (1704 .. 2104). The "Valor Mora" (F) value travels together with
# its period, so the whole block of rows is effectively re-sorted ascending
# by period (Salario Basico in column G is unchanged for every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order of periods (previously descending 2104 -> 1704)
$periods = @(
    "1704","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906",
    "1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005",
    "2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104"
)

# Valor Mora for each period in the same new order (most periods carry the
# same mora value of 31249; only the first/last periods differ)
$valorMora = @(
    29509,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,27083
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
}
